$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (fill/border/style) from the last data row (25) down into
# the new row (26), matching the "Övrigt" category block styling, then fill
# in the new Kategori/Underkategori/Aktivitet values for the added "Rast"
# activity.
$ws.Range("A25:C25").Copy() | Out-Null
$ws.Range("A26:C26").PasteSpecial(-4122) | Out-Null

$ws.Range("A26").Value = "Övrigt"
$ws.Range("B26").Value = "Övrigt"
$ws.Range("C26").Value = "Rast"

# Grow Table1 so the new row is included in the table (ref + autofilter
# both expand from A1:C25 to A1:C26).
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:C26")) | Out-Null

# Move the active selection to G12 (was E12).
$ws.Range("G12").Select() | Out-Null
